$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix label text: "Extension of remarks" -> "Extensions of remarks"
$ws.Range("A9").Value = "Extensions of remarks"

# Fix "Time in session" strings (spacing/punctuation correction)
$ws.Range("C7").Value = "1247 hrs., 52'"
$ws.Range("B7").Value = "1375 hrs., 54'"

# Convert "Pages of proceedings" from text-with-letter-prefix to plain numbers
$ws.Range("B8").Value = 16071
$ws.Range("C8").Value = 16951

# Convert "Extensions of remarks" House value from text-with-letter-prefix to plain number
$ws.Range("C9").Value = 2664

$wb.Save()
